# Daily attendance processing - reorder "Recorded By" (column G) entries
# so that entries containing "System" as one of the comma-separated values
# have their order reversed (e.g. "x, System" -> "System, x").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = 7
    $value = $cell.Value()

    if ($value -ne $null -and $value -ne "") {
        $parts = $value -split ", "
        if ($parts.Count -gt 1) {
            $hasSystem = $false
            foreach ($p in $parts) {
                if ($p -eq "System") {
                    $hasSystem = $true
                }
            }
            if ($hasSystem) {
                $reversedParts = $parts[-1..-$parts.Count]
                $newValue = [string]::Join(", ", $reversedParts)
                $cell.Value = $newValue
            }
        }
    }
}
